$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended to the sheet (rows 114-117)
# Columns: A B C D E F G H I  (C and D are left blank, like the rows above them)
$newRows = @(
    @(649, 329, $null, $null, 6, 113, 314, 38, 0),
    @(649, 329, $null, $null, 6, 114, 314, 0, 0),
    @(649, 329, $null, $null, 6, 114, 314, 0, 0),
    @(649, 329, $null, $null, 6, 116, 314, 0, 0)
)

$startRow = 114
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $val = $rowData[$c]
        if ($null -eq $val) {
            # Materialize an empty cell (matching existing blank C/D cells) without
            # introducing a new style entry.
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
